$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Archive ID" -> "LED Set"
$ws.Range("A1").Value = "LED Set"

# Convert Channel column (C) from text Top/Bot to numeric 1/2 for all data rows
for ($r = 2; $r -le 43; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -eq "Top") {
        $cell.Value = 1
    } elseif ($current -eq "Bot") {
        $cell.Value = 2
    }
}

# Update the selected cell in the sheet view
[void]$ws.Range("E14").Select()
